$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7.818245253082684
$ws.Range("D2").Value = 5.031896916800619
$ws.Range("E2").Value = 12.31065508756321
$ws.Range("F2").Value = 27.98736913050699
$ws.Range("G2").Value = 35.73474102701444
$ws.Range("H2").Value = 15.5529244094609
$ws.Range("K2").Value = 9.984160097481048
$ws.Range("M2").Value = 13.83689475540853
$ws.Range("B3").Value = 7.74815202461363
$ws.Range("D3").Value = 5.036156397470084
$ws.Range("E3").Value = 12.16266361848166
$ws.Range("F3").Value = 27.81766522492601
$ws.Range("G3").Value = 35.37028824802812
$ws.Range("H3").Value = 15.55865381647555
$ws.Range("K3").Value = 9.508441533254315
$ws.Range("M3").Value = 13.61121181490203
$ws.Range("B4").Value = 7.706697473712548
$ws.Range("D4").Value = 5.038837792927202
$ws.Range("E4").Value = 12.07503423584158
$ws.Range("F4").Value = 27.72161652226385
$ws.Range("G4").Value = 35.15747961846768
$ws.Range("H4").Value = 15.56559795591642
$ws.Range("K4").Value = 9.201629725911936
$ws.Range("M4").Value = 13.47466976403396
$ws.Range("B5").Value = 7.690220682354167
$ws.Range("D5").Value = 5.039947387510532
$ws.Range("E5").Value = 12.04017882029557
$ws.Range("F5").Value = 27.68455555728748
$ws.Range("G5").Value = 35.0736077772185
$ws.Range("H5").Value = 15.56928761792437
$ws.Range("K5").Value = 9.072976542764602
$ws.Range("M5").Value = 13.41961326984506
$ws.Range("B6").Value = 7.687510382666407
$ws.Range("D6").Value = 5.04013266451944
$ws.Range("E6").Value = 12.03444382900811
$ws.Range("F6").Value = 27.67852802584947
$ws.Range("G6").Value = 35.05985555779374
$ws.Range("H6").Value = 15.56995215996632
$ws.Range("K6").Value = 9.051397350373733
$ws.Range("M6").Value = 13.41050873622424
$ws.Range("B7").Value = 7.706473552966937
$ws.Range("D7").Value = 5.038852688492468
$ws.Range("E7").Value = 12.07456065328472
$ws.Range("F7").Value = 27.7211082477106
$ws.Range("G7").Value = 35.15633684063823
$ws.Range("H7").Value = 15.56564423701555
$ws.Range("K7").Value = 9.19990923008978
$ws.Range("M7").Value = 13.47392478300221
$ws.Range("B8").Value = 7.793759728250695
$ws.Range("D8").Value = 5.033352050707428
$ws.Range("E8").Value = 12.25898025077999
$ws.Range("F8").Value = 27.9271829604257
$ws.Range("G8").Value = 35.60685869550485
$ws.Range("H8").Value = 15.55418788820845
$ws.Range("K8").Value = 9.823228764579008
$ws.Range("M8").Value = 13.75870377931085
$ws.Range("B9").Value = 7.976651945707117
$ws.Range("D9").Value = 5.023076411470131
$ws.Range("E9").Value = 12.6444633931564
$ws.Range("F9").Value = 28.3944538236407
$ws.Range("G9").Value = 36.57298883514213
$ws.Range("H9").Value = 15.55896756377476
$ws.Range("K9").Value = 10.92607534549652
$ws.Range("M9").Value = 14.32991751858079
$ws.Range("B10").Value = 8.117009746435025
$ws.Range("D10").Value = 5.015820722131439
$ws.Range("E10").Value = 12.93962169660994
$ws.Range("F10").Value = 28.77399751765469
$ws.Range("G10").Value = 37.32696098465087
$ws.Range("H10").Value = 15.57914692950957
$ws.Range("K10").Value = 11.66066493551621
$ws.Range("M10").Value = 14.75294450290721
$ws.Range("B11").Value = 8.181907028305101
$ws.Range("D11").Value = 5.012580038135287
$ws.Range("E11").Value = 13.07592171397211
$ws.Range("F11").Value = 28.95398280560413
$ws.Range("G11").Value = 37.67813674334544
$ws.Range("H11").Value = 15.59194965757844
$ws.Range("K11").Value = 11.97800083859597
$ws.Range("M11").Value = 14.9452076405874
$ws.Range("B12").Value = 8.206610941997361
$ws.Range("D12").Value = 5.011361217352515
$ws.Range("E12").Value = 13.12777903206715
$ws.Range("F12").Value = 29.02314231070861
$ws.Range("G12").Value = 37.81217331062211
$ws.Range("H12").Value = 15.5973181304582
$ws.Range("K12").Value = 12.09572143558927
$ws.Range("M12").Value = 15.01791454169062
$ws.Range("B13").Value = 8.201285136809673
$ws.Range("D13").Value = 5.011623344762328
$ws.Range("E13").Value = 13.1166005737204
$ws.Range("F13").Value = 29.00820380015496
$ws.Range("G13").Value = 37.78326124757967
$ws.Range("H13").Value = 15.59613880578654
$ws.Range("K13").Value = 12.0704774265197
$ws.Range("M13").Value = 15.00226135934925
$ws.Range("B14").Value = 8.183936980144026
$ws.Range("D14").Value = 5.012479599059076
$ws.Range("E14").Value = 13.08018344487756
$ws.Range("F14").Value = 28.95965279799664
$ws.Range("G14").Value = 37.6891436035237
$ws.Range("H14").Value = 15.59238090162915
$ws.Range("K14").Value = 11.98773497123153
$ws.Range("M14").Value = 14.9511916784105
$ws.Range("B15").Value = 8.173326844889081
$ws.Range("D15").Value = 5.013005159905734
$ws.Range("E15").Value = 13.05790716945637
$ws.Range("F15").Value = 28.93004297588534
$ws.Range("G15").Value = 37.63162738596682
$ws.Range("H15").Value = 15.59014681594693
$ws.Range("K15").Value = 11.93673337244077
$ws.Range("M15").Value = 14.91989498708612
$ws.Range("B16").Value = 8.112787765215788
$ws.Range("D16").Value = 5.016033692135622
$ws.Range("E16").Value = 12.93075088043138
$ws.Range("F16").Value = 28.76237823351939
$ws.Range("G16").Value = 37.30416507799252
$ws.Range("H16").Value = 15.57838310276156
$ws.Range("K16").Value = 11.63958511655185
$ws.Range("M16").Value = 14.74036976931827
$ws.Range("B17").Value = 8.07590176568058
$ws.Range("D17").Value = 5.017906765165058
$ws.Range("E17").Value = 12.85323018245942
$ws.Range("F17").Value = 28.66136298609518
$ws.Range("G17").Value = 37.1052876360984
$ws.Range("H17").Value = 15.57209409360707
$ws.Range("K17").Value = 11.45296180292278
$ws.Range("M17").Value = 14.63013883036211
$ws.Range("B18").Value = 8.054785859549614
$ws.Range("D18").Value = 5.01898976694643
$ws.Range("E18").Value = 12.80883640584042
$ws.Range("F18").Value = 28.60395472914573
$ws.Range("G18").Value = 36.99167701107541
$ws.Range("H18").Value = 15.5688179064823
$ws.Range("K18").Value = 11.34403784546952
$ws.Range("M18").Value = 14.56672480073775
$ws.Range("B19").Value = 8.047654210524126
$ws.Range("D19").Value = 5.019357433068032
$ws.Range("E19").Value = 12.79384025269112
$ws.Range("F19").Value = 28.58463782365714
$ws.Range("G19").Value = 36.95334784758624
$ws.Range("H19").Value = 15.56776723893254
$ws.Range("K19").Value = 11.30688708171726
$ws.Range("M19").Value = 14.54525411821037
$ws.Range("B20").Value = 8.079818158297064
$ws.Range("D20").Value = 5.017706789516736
$ws.Range("E20").Value = 12.86146265944116
$ws.Range("F20").Value = 28.67204486588209
$ws.Range("G20").Value = 37.12637880806864
$ws.Range("H20").Value = 15.57272827139425
$ws.Range("K20").Value = 11.47299224024666
$ws.Range("M20").Value = 14.64187489972968
$ws.Range("B21").Value = 8.18902923579842
$ws.Range("D21").Value = 5.012227871734869
$ws.Range("E21").Value = 13.0908738169614
$ws.Range("F21").Value = 28.97388660035963
$ws.Range("G21").Value = 37.71676062632937
$ws.Range("H21").Value = 15.59347057418594
$ws.Range("K21").Value = 12.01210505056819
$ws.Range("M21").Value = 14.96619534337682
$ws.Range("B22").Value = 8.261145115345975
$ws.Range("D22").Value = 5.008695667452159
$ws.Range("E22").Value = 13.24220405805628
$ws.Range("F22").Value = 29.17697866212173
$ws.Range("G22").Value = 38.10869430044895
$ws.Range("H22").Value = 15.61005928346178
$ws.Range("K22").Value = 12.35017123984631
$ws.Range("M22").Value = 15.1775483355985
$ws.Range("B23").Value = 8.222595039514948
$ws.Range("D23").Value = 5.010576512436043
$ws.Range("E23").Value = 13.16132442473805
$ws.Range("F23").Value = 29.06806935965002
$ws.Range("G23").Value = 37.89899630903558
$ws.Range("H23").Value = 15.60092844170985
$ws.Range("K23").Value = 12.17105260093652
$ws.Range("M23").Value = 15.06482469854759
$ws.Range("B24").Value = 8.078047274114395
$ws.Range("D24").Value = 5.017797179298073
$ws.Range("E24").Value = 12.8577402112399
$ws.Range("F24").Value = 28.66721350715265
$ws.Range("G24").Value = 37.11684121615573
$ws.Range("H24").Value = 15.57244050211574
$ws.Range("K24").Value = 11.46394156002824
$ws.Range("M24").Value = 14.63656914416285
$ws.Range("B25").Value = 7.926037600269058
$ws.Range("D25").Value = 5.025803387706343
$ws.Range("E25").Value = 12.53789100268169
$ws.Range("F25").Value = 28.26150369359329
$ws.Range("G25").Value = 36.30340109044072
$ws.Range("H25").Value = 15.55474994303191
$ws.Range("K25").Value = 10.64084721315802
$ws.Range("M25").Value = 14.17448379281098
